$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 2449.4443
$ws.Range("I19").Value = 2816.8333
$ws.Range("J19").Value = 1714.6666
$ws.Range("K19").Value = 2816.8333
$ws.Range("L19").Value = 1714.6666
$ws.Range("M19").Value = -2641.8333
$ws.Range("N19").Value = -2064.6666
$ws.Range("H39").Value = 148.33333
$ws.Range("I39").Value = 148.33333
$ws.Range("K39").Value = 444.99999
$ws.Range("M39").Value = -148.99999
$ws.Range("H42").Value = 555.9
$ws.Range("J42").Value = 665.5714
$ws.Range("L42").Value = 1996.7142
$ws.Range("N42").Value = -2456.7142
$ws.Range("H43").Value = 7249
$ws.Range("I43").Value = 6998
$ws.Range("K43").Value = 6998
$ws.Range("M43").Value = -6929
$ws.Range("H64").Value = 2767.5334
$ws.Range("I64").Value = 2500.5
$ws.Range("K64").Value = 2500.5
$ws.Range("M64").Value = -2252.5
$ws.Range("H67").Value = 2767.5334
$ws.Range("I67").Value = 2500.5
$ws.Range("K67").Value = 2500.5
$ws.Range("M67").Value = -1642.5
$ws.Range("H92").Value = 1010.25
$ws.Range("I92").Value = 964.1667
$ws.Range("K92").Value = 964.1667
$ws.Range("M92").Value = 283.8333
$ws.Range("H112").Value = 1724.2354
$ws.Range("J112").Value = 1922.2858
$ws.Range("L112").Value = 5766.857400000001
$ws.Range("N112").Value = -7982.857400000001
$ws.Range("H132").Value = 1085.3214
$ws.Range("I132").Value = 1085.3214
$ws.Range("K132").Value = 3255.9642
$ws.Range("M132").Value = -725.9642000000003
$ws.Range("H138").Value = 3289.4138
$ws.Range("J138").Value = 3380.2744
$ws.Range("L138").Value = 10140.8232
$ws.Range("N138").Value = -20420.8232

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1257.25
$ws.Range("I2").Value = 1210.3
$ws.Range("J2").Value = 1492
$ws.Range("K2").Value = 1210.3
$ws.Range("L2").Value = 1492
$ws.Range("M2").Value = -1097.3
$ws.Range("N2").Value = -1718
$ws.Range("H32").Value = 4726.18
$ws.Range("I32").Value = 3419.691
$ws.Range("K32").Value = 3419.691
$ws.Range("M32").Value = -3132.691
$ws.Range("H45").Value = 2428.1667
$ws.Range("I45").Value = 2428.1667
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 2428.1667
$ws.Range("L45").Value = 0
$ws.Range("M45").ClearContents()
$ws.Range("N45").Value = -2051.1667
$ws.Range("H63").Value = 1147.75
$ws.Range("I63").Value = 1025.7142
$ws.Range("K63").Value = 1025.7142
$ws.Range("M63").Value = -339.7141999999999
$ws.Range("H66").Value = 1147.75
$ws.Range("I66").Value = 1025.7142
$ws.Range("K66").Value = 5128.571
$ws.Range("M66").Value = -1696.571
$ws.Range("H74").Value = 1226.3334
$ws.Range("I74").Value = 967.125
$ws.Range("K74").Value = 967.125
$ws.Range("M74").Value = -93.125
$ws.Range("H77").Value = 1226.3334
$ws.Range("I77").Value = 967.125
$ws.Range("K77").Value = 4835.625
$ws.Range("M77").Value = -467.625
$ws.Range("H102").Value = 2935.25
$ws.Range("I102").Value = 2935.25
$ws.Range("K102").Value = 2935.25
$ws.Range("M102").Value = -1313.25
$ws.Range("H110").Value = 1119.4
$ws.Range("I110").Value = 1066
$ws.Range("K110").Value = 1066
$ws.Range("M110").Value = 979
$ws.Range("H116").Value = 1257.25
$ws.Range("I116").Value = 1210.3
$ws.Range("J116").Value = 1492
$ws.Range("K116").Value = 1210.3
$ws.Range("L116").Value = 1492
$ws.Range("M116").Value = 1083.7
$ws.Range("N116").Value = -6080
$ws.Range("H122").Value = 1028.5
$ws.Range("I122").Value = 798.4286
$ws.Range("J122").Value = 1565.3334
$ws.Range("K122").Value = 2395.2858
$ws.Range("L122").Value = 4696.0002
$ws.Range("M122").Value = 54.71420000000035
$ws.Range("N122").Value = -9596.0002

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1257.25
$ws.Range("I3").Value = 1210.3
$ws.Range("J3").Value = 1492
$ws.Range("K3").Value = 1210.3
$ws.Range("L3").Value = 1492
$ws.Range("M3").Value = -1096.3
$ws.Range("N3").Value = -1720
$ws.Range("H64").Value = 1246.5
$ws.Range("I64").Value = 1003
$ws.Range("K64").Value = 1003
$ws.Range("M64").Value = -778
$ws.Range("H67").Value = 1246.5
$ws.Range("I67").Value = 1003
$ws.Range("K67").Value = 1003
$ws.Range("M67").Value = -223
$ws.Range("H94").Value = 1177.1666
$ws.Range("I94").Value = 746
$ws.Range("K94").Value = 746
$ws.Range("M94").Value = -295
$ws.Range("H107").Value = 885
$ws.Range("I107").Value = 885
$ws.Range("K107").Value = 885
$ws.Range("M107").Value = 1035
$ws.Range("H138").Value = 124849.5
$ws.Range("J138").Value = 124849.5
$ws.Range("L138").Value = 124849.5
$ws.Range("N138").Value = -135129.5

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 680
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("H55").Value = 24992
$ws.Range("J55").Value = 24992
$ws.Range("L55").Value = 24992
$ws.Range("N55").Value = -25622
$ws.Range("H58").Value = 2563.8
$ws.Range("I58").Value = 2159.2
$ws.Range("J58").Value = 2968.4
$ws.Range("K58").Value = 2159.2
$ws.Range("L58").Value = 2968.4
$ws.Range("M58").Value = -1956.2
$ws.Range("N58").Value = -3374.4
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("H136").Value = 2563.8
$ws.Range("I136").Value = 2159.2
$ws.Range("J136").Value = 2968.4
$ws.Range("K136").Value = 6477.599999999999
$ws.Range("L136").Value = 8905.200000000001
$ws.Range("M136").Value = -3927.599999999999
$ws.Range("N136").Value = -14005.2

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 6106.375
$ws.Range("I80").Value = 4491.8
$ws.Range("K80").Value = 4491.8
$ws.Range("M80").Value = -3493.8
$ws.Range("H83").Value = 6106.375
$ws.Range("I83").Value = 4491.8
$ws.Range("K83").Value = 22459
$ws.Range("M83").Value = -17467
$ws.Range("H102").Value = 1949.7778
$ws.Range("I102").Value = 1649.8572
$ws.Range("K102").Value = 1649.8572
$ws.Range("M102").Value = -27.85719999999992
$ws.Range("H113").Value = 1586.875
$ws.Range("I113").Value = 1586.875
$ws.Range("K113").Value = 1586.875
$ws.Range("M113").Value = 583.125
$ws.Range("H122").Value = 1601.5834
$ws.Range("I122").Value = 1277.625
$ws.Range("J122").Value = 2249.5
$ws.Range("K122").Value = 3832.875
$ws.Range("L122").Value = 6748.5
$ws.Range("M122").Value = -1382.875
$ws.Range("N122").Value = -11648.5
$ws.Range("H126").Value = 998.25
$ws.Range("I126").Value = 997
$ws.Range("K126").Value = 2991
$ws.Range("M126").Value = -521

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1250
$ws.Range("I7").Value = 1250
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 1250
$ws.Range("L7").Value = 0
$ws.Range("M7").ClearContents()
$ws.Range("N7").Value = -1138
$ws.Range("H22").Value = 4562.909
$ws.Range("I22").Value = 2500
$ws.Range("K22").Value = 2500
$ws.Range("M22").Value = -2205
$ws.Range("H27").Value = 4562.909
$ws.Range("I27").Value = 2500
$ws.Range("K27").Value = 2500
$ws.Range("M27").Value = -2393
$ws.Range("H43").Value = 954998.9
$ws.Range("J43").Value = 954998.9
$ws.Range("L43").Value = 954998.9
$ws.Range("N43").Value = -955384.9
$ws.Range("H126").Value = 1250
$ws.Range("I126").Value = 1250
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 3750
$ws.Range("L126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = -1280

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 9018.444
$ws.Range("I81").Value = 4030.1667
$ws.Range("K81").Value = 8060.3334
$ws.Range("M81").Value = -6999.3334
$ws.Range("H84").Value = 9018.444
$ws.Range("I84").Value = 4030.1667
$ws.Range("K84").Value = 40301.667
$ws.Range("M84").Value = -34997.667
$ws.Range("H126").Value = 1912.6875
$ws.Range("I126").Value = 1758.75
$ws.Range("K126").Value = 5276.25
$ws.Range("M126").Value = -2806.25
$ws.Range("H132").Value = 3132.25
$ws.Range("I132").Value = 2758.9
$ws.Range("K132").Value = 8276.700000000001
$ws.Range("M132").Value = -5746.700000000001
$ws.Range("H136").Value = 2408.24
$ws.Range("I136").Value = 2111.95
$ws.Range("J136").Value = 3593.4
$ws.Range("K136").Value = 6335.849999999999
$ws.Range("L136").Value = 10780.2
$ws.Range("M136").Value = -3785.849999999999
$ws.Range("N136").Value = -15880.2
